$wb = $excel.ActiveWorkbook

# Work on the Orders sheet
$ws = $wb.Worksheets.Item("Orders")
$ws.Activate()

# Row 9: add Status/Notes to the existing "Ipoh Coff" row
$ws.Range("C9").Value = "Failed"
$ws.Range("D9").Value = "Product Does Not Exist"

# Row 10: new row - Chai, out of stock quantity
$ws.Range("A10").Value = "Chai"
$ws.Range("B10").Value = 3000
$ws.Range("C10").Value = "Failed"
$ws.Range("D10").Value = "Out of Stock"

# Row 11: new row - Laughing Lumberjack Lager, out of stock quantity
$ws.Range("A11").Value = "Laughing Lumberjack Lager"
$ws.Range("B11").Value = 5000
$ws.Range("C11").Value = "Failed"
$ws.Range("D11").Value = "Out of Stock"

$ws.Range("G9").Select()
